$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/value updates (safe as-is, not numeric-looking or percentage strings)
$ws.Range("D2").Value = "66.305.03"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.223.33"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +4.35%  "
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.221.81"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").Value = "3.754.32"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("E16").Value = "  +4.34%  "
$ws.Range("D17").Value = "66.359.50"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "3.222.52"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  +6.13%  "
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +3.72%  "
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("E29").Value = "  +2.69%  "
$ws.Range("E30").Value = "  +11.15%  "
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.986.11"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E45").Value = "  +4.39%  "
$ws.Range("D46").Value = "0.0₃0650"
$ws.Range("E46").Value = "  +7.22%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("E51").Value = "  +4.90%  "

# Numeric-looking price text must be forced to remain as text (matching original inlineStr type)
# Use NumberFormat "@" to force text entry, then ClearFormats to drop the residual style index
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.32"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.90"
$ws.Range("D6").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.18"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.511"
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.01"
$ws.Range("D14").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.50"
$ws.Range("D16").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "514.23"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.82"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.28"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.02"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.60"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.31"
$ws.Range("D28").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.94"
$ws.Range("D30").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.30"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.22"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.67"
$ws.Range("D35").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0923"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "488.86"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0423"
$ws.Range("D39").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.88"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.297"
$ws.Range("D43").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.27"
$ws.Range("D47").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.93"
$ws.Range("D51").ClearFormats()
